# Update the header row labels to remove leading spaces and fix "Market Cap" column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Market Cap"
$ws.Range("G1").Value = "Trade Date"
$ws.Range("F1").Value = "Volume"
$ws.Range("E1").Value = "Close"
$ws.Range("D1").Value = "Open"
$ws.Range("C1").Value = "Sector"
$ws.Range("B1").Value = "Company Name"

# Convert the Market Cap column from text (e.g. " 231.8B") to numeric values
$ws.Range("H2").Value = 231.8
$ws.Range("H3").Value = 175.2
$ws.Range("H4").Value = 85.3
$ws.Range("H5").Value = 141.5
$ws.Range("H6").Value = 47.2
$ws.Range("H7").Value = 82.1

# Update the selected cell to match the target state
$ws.Range("E6").Select()
